$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 14.15440837047039
$ws.Cells.Item(2, 3).Value = 4.763913061872528
$ws.Cells.Item(2, 4).Value = 9.106701891441801
$ws.Cells.Item(2, 5).Value = 10.42840865307176
$ws.Cells.Item(2, 6).Value = 45.31107084622924
$ws.Cells.Item(2, 8).Value = 7.344005520526261
$ws.Cells.Item(2, 9).Value = 37.03918245434775
$ws.Cells.Item(2, 11).Value = 12.59491124901308
$ws.Cells.Item(2, 12).Value = 10.54869075530686
$ws.Cells.Item(2, 13).Value = 15.46849324139407
$ws.Cells.Item(3, 2).Value = 14.07919227042468
$ws.Cells.Item(3, 3).Value = 4.521650149323654
$ws.Cells.Item(3, 4).Value = 9.096251255058386
$ws.Cells.Item(3, 5).Value = 10.44012375205418
$ws.Cells.Item(3, 6).Value = 44.98296002507969
$ws.Cells.Item(3, 8).Value = 7.344005520526261
$ws.Cells.Item(3, 9).Value = 36.87890517410655
$ws.Cells.Item(3, 11).Value = 12.5384462716342
$ws.Cells.Item(3, 12).Value = 10.55690322042603
$ws.Cells.Item(3, 13).Value = 15.48061640937515
$ws.Cells.Item(4, 2).Value = 14.03759974515406
$ws.Cells.Item(4, 3).Value = 4.36505995200638
$ws.Cells.Item(4, 4).Value = 9.089571111249564
$ws.Cells.Item(4, 5).Value = 10.44809934928624
$ws.Cells.Item(4, 6).Value = 44.78664169360854
$ws.Cells.Item(4, 8).Value = 7.344005520526261
$ws.Cells.Item(4, 9).Value = 36.78399229748975
$ws.Cells.Item(4, 11).Value = 12.50763355593799
$ws.Cells.Item(4, 12).Value = 10.56341049827725
$ws.Cells.Item(4, 13).Value = 15.49112978020088
$ws.Cells.Item(5, 2).Value = 14.02182225988627
$ws.Cells.Item(5, 3).Value = 4.299305600786679
$ws.Cells.Item(5, 4).Value = 9.086782373268074
$ws.Cells.Item(5, 5).Value = 10.45154654763899
$ws.Cells.Item(5, 6).Value = 44.70798604298045
$ws.Cells.Item(5, 8).Value = 7.344005520526261
$ws.Cells.Item(5, 9).Value = 36.74621056018894
$ws.Cells.Item(5, 11).Value = 12.49605880608902
$ws.Cells.Item(5, 12).Value = 10.56643086163656
$ws.Cells.Item(5, 13).Value = 15.49618636401391
$ws.Cells.Item(6, 2).Value = 14.01927366673519
$ws.Cells.Item(6, 3).Value = 4.288270775519917
$ws.Cells.Item(6, 4).Value = 9.086315278974954
$ws.Cells.Item(6, 5).Value = 10.45213086453234
$ws.Cells.Item(6, 6).Value = 44.69500808917584
$ws.Cells.Item(6, 8).Value = 7.344005520526261
$ws.Cells.Item(6, 9).Value = 36.73999150772359
$ws.Cells.Item(6, 11).Value = 12.49419643437504
$ws.Cells.Item(6, 12).Value = 10.56695465912585
$ws.Cells.Item(6, 13).Value = 15.49707265752916
$ws.Cells.Item(7, 2).Value = 14.03738219826458
$ws.Cells.Item(7, 3).Value = 4.364180988735096
$ws.Cells.Item(7, 4).Value = 9.089533771106293
$ws.Cells.Item(7, 5).Value = 10.44814504102388
$ws.Cells.Item(7, 6).Value = 44.78557540022118
$ws.Cells.Item(7, 8).Value = 7.344005520526261
$ws.Cells.Item(7, 9).Value = 36.78347911100355
$ws.Cells.Item(7, 11).Value = 12.50747346556249
$ws.Cells.Item(7, 12).Value = 10.56344973923695
$ws.Cells.Item(7, 13).Value = 15.49119484772513
$ws.Cells.Item(8, 2).Value = 14.12753221161241
$ws.Cells.Item(8, 3).Value = 4.682023433956594
$ws.Cells.Item(8, 4).Value = 9.10315300927307
$ws.Cells.Item(8, 5).Value = 10.43228580670822
$ws.Cells.Item(8, 6).Value = 45.1969029562372
$ws.Cells.Item(8, 8).Value = 7.344005520526261
$ws.Cells.Item(8, 9).Value = 36.98320299637204
$ws.Cells.Item(8, 11).Value = 12.57464913939016
$ws.Cells.Item(8, 12).Value = 10.55121857141885
$ws.Cells.Item(8, 13).Value = 15.4720365066645
$ws.Cells.Item(9, 2).Value = 14.33987600051752
$ws.Cells.Item(9, 3).Value = 5.242122731514795
$ws.Cells.Item(9, 4).Value = 9.127783482380348
$ws.Cells.Item(9, 5).Value = 10.40738069830952
$ws.Cells.Item(9, 6).Value = 46.04191804463387
$ws.Cells.Item(9, 8).Value = 7.344005520526261
$ws.Cells.Item(9, 9).Value = 37.40185415276833
$ws.Cells.Item(9, 11).Value = 12.73639499215219
$ws.Cells.Item(9, 12).Value = 10.53884257471795
$ws.Cells.Item(9, 13).Value = 15.45879720280734
$ws.Cells.Item(10, 2).Value = 14.51627967081294
$ws.Cells.Item(10, 3).Value = 5.614117291717593
$ws.Cells.Item(10, 4).Value = 9.144632050779562
$ws.Cells.Item(10, 5).Value = 10.39284100955307
$ws.Cells.Item(10, 6).Value = 46.68281271699482
$ws.Cells.Item(10, 8).Value = 7.344005520526261
$ws.Cells.Item(10, 9).Value = 37.7248167794261
$ws.Cells.Item(10, 11).Value = 12.87264201442705
$ws.Cells.Item(10, 12).Value = 10.53680699720717
$ws.Cells.Item(10, 13).Value = 15.46385355267193
$ws.Cells.Item(11, 2).Value = 14.60065147037849
$ws.Cells.Item(11, 3).Value = 5.774634069343143
$ws.Cells.Item(11, 4).Value = 9.152029977660311
$ws.Cells.Item(11, 5).Value = 10.3870385955086
$ws.Cells.Item(11, 6).Value = 46.97798147035312
$ws.Cells.Item(11, 8).Value = 7.344005520526261
$ws.Cells.Item(11, 9).Value = 37.87483179524321
$ws.Cells.Item(11, 11).Value = 12.93818980228987
$ws.Cells.Item(11, 12).Value = 10.53740751975719
$ws.Cells.Item(11, 13).Value = 15.46934778368548
$ws.Cells.Item(12, 2).Value = 14.63316619133229
$ws.Cells.Item(12, 3).Value = 5.834158085481824
$ws.Cells.Item(12, 4).Value = 9.154793340403252
$ws.Cells.Item(12, 5).Value = 10.38495776642915
$ws.Cells.Item(12, 6).Value = 47.09021013130783
$ws.Cells.Item(12, 8).Value = 7.344005520526261
$ws.Cells.Item(12, 9).Value = 37.93206140281046
$ws.Cells.Item(12, 11).Value = 12.96350435470621
$ws.Cells.Item(12, 12).Value = 10.53785377454323
$ws.Cells.Item(12, 13).Value = 15.47188572228558
$ws.Cells.Item(13, 2).Value = 14.62613889899466
$ws.Cells.Item(13, 3).Value = 5.821394680119312
$ws.Cells.Item(13, 4).Value = 9.154199889919902
$ws.Cells.Item(13, 5).Value = 10.3854007374449
$ws.Cells.Item(13, 6).Value = 47.06602052323527
$ws.Cells.Item(13, 8).Value = 7.344005520526261
$ws.Cells.Item(13, 9).Value = 37.91971758506379
$ws.Cells.Item(13, 11).Value = 12.95803081562091
$ws.Cells.Item(13, 12).Value = 10.53774794361831
$ws.Cells.Item(13, 13).Value = 15.4713188227276
$ws.Cells.Item(14, 2).Value = 14.60331531911257
$ws.Cells.Item(14, 3).Value = 5.77955644199828
$ws.Cells.Item(14, 4).Value = 9.152258082706968
$ws.Cells.Item(14, 5).Value = 10.38686507303639
$ws.Cells.Item(14, 6).Value = 46.98720586009048
$ws.Cells.Item(14, 8).Value = 7.344005520526261
$ws.Cells.Item(14, 9).Value = 37.87953178447223
$ws.Cells.Item(14, 11).Value = 12.94026268330019
$ws.Cells.Item(14, 12).Value = 10.53743985141926
$ws.Cells.Item(14, 13).Value = 15.46954742659592
$ws.Cells.Item(15, 2).Value = 14.58940793441239
$ws.Cells.Item(15, 3).Value = 5.753764988956685
$ws.Cells.Item(15, 4).Value = 9.151063719353642
$ws.Cells.Item(15, 5).Value = 10.38777717275215
$ws.Cells.Item(15, 6).Value = 46.938986833215
$ws.Cells.Item(15, 8).Value = 7.344005520526261
$ws.Cells.Item(15, 9).Value = 37.85497109734834
$ws.Cells.Item(15, 11).Value = 12.92944278289144
$ws.Cells.Item(15, 12).Value = 10.53727961569181
$ws.Cells.Item(15, 13).Value = 15.46852189954252
$ws.Cells.Item(16, 2).Value = 14.51084639430343
$ws.Cells.Item(16, 3).Value = 5.603451196724925
$ws.Cells.Item(16, 4).Value = 9.144143237161503
$ws.Cells.Item(16, 5).Value = 10.39323651921307
$ws.Cells.Item(16, 6).Value = 46.66358999514382
$ws.Cells.Item(16, 8).Value = 7.344005520526261
$ws.Cells.Item(16, 9).Value = 37.71507330785395
$ws.Cells.Item(16, 11).Value = 12.86842842138641
$ws.Cells.Item(16, 12).Value = 10.53679841447859
$ws.Cells.Item(16, 13).Value = 15.46355860806073
$ws.Cells.Item(17, 2).Value = 14.46368822684977
$ws.Cells.Item(17, 3).Value = 5.509002152256529
$ws.Cells.Item(17, 4).Value = 9.139829657106873
$ws.Cells.Item(17, 5).Value = 10.39679334360056
$ws.Cells.Item(17, 6).Value = 46.49552391494316
$ws.Cells.Item(17, 8).Value = 7.344005520526261
$ws.Cells.Item(17, 9).Value = 37.63002763727706
$ws.Cells.Item(17, 11).Value = 12.83189804694174
$ws.Cells.Item(17, 12).Value = 10.53689374714501
$ws.Cells.Item(17, 13).Value = 15.46133060240835
$ws.Cells.Item(18, 2).Value = 14.43695411838073
$ws.Cells.Item(18, 3).Value = 5.453858497625558
$ws.Cells.Item(18, 4).Value = 9.137323493620791
$ws.Cells.Item(18, 5).Value = 10.39891556088523
$ws.Cells.Item(18, 6).Value = 46.39920175177483
$ws.Cells.Item(18, 8).Value = 7.344005520526261
$ws.Cells.Item(18, 9).Value = 37.58140463179206
$ws.Cells.Item(18, 11).Value = 12.81122388059484
$ws.Cells.Item(18, 12).Value = 10.53709230795095
$ws.Cells.Item(18, 13).Value = 15.46034988691187
$ws.Cells.Item(19, 2).Value = 14.42797023529305
$ws.Cells.Item(19, 3).Value = 5.43504742390701
$ws.Cells.Item(19, 4).Value = 9.136470629673051
$ws.Cells.Item(19, 5).Value = 10.39964724282485
$ws.Cells.Item(19, 6).Value = 46.36664998123803
$ws.Cells.Item(19, 8).Value = 7.344005520526261
$ws.Cells.Item(19, 9).Value = 37.56499273912033
$ws.Cells.Item(19, 11).Value = 12.80428246278381
$ws.Cells.Item(19, 12).Value = 10.53718424139312
$ws.Cells.Item(19, 13).Value = 15.46006954513181
$ws.Cells.Item(20, 2).Value = 14.46866812916586
$ws.Cells.Item(20, 3).Value = 5.519141252245399
$ws.Cells.Item(20, 4).Value = 9.140291441583885
$ws.Cells.Item(20, 5).Value = 10.39640680601326
$ws.Cells.Item(20, 6).Value = 46.51337960733576
$ws.Cells.Item(20, 8).Value = 7.344005520526261
$ws.Cells.Item(20, 9).Value = 37.63905073439264
$ws.Cells.Item(20, 11).Value = 12.83575201754993
$ws.Cells.Item(20, 12).Value = 10.53686872805998
$ws.Cells.Item(20, 13).Value = 15.46153666136152
$ws.Cells.Item(21, 2).Value = 14.61000405549268
$ws.Cells.Item(21, 3).Value = 5.791879603431142
$ws.Cells.Item(21, 4).Value = 9.152829469746205
$ws.Cells.Item(21, 5).Value = 10.3864318051019
$ws.Cells.Item(21, 6).Value = 47.01034382974478
$ws.Cells.Item(21, 8).Value = 7.344005520526261
$ws.Cells.Item(21, 9).Value = 37.89132405971952
$ws.Cells.Item(21, 11).Value = 12.94546839880322
$ws.Cells.Item(21, 12).Value = 10.53752441173493
$ws.Cells.Item(21, 13).Value = 15.47005533177802
$ws.Cells.Item(22, 2).Value = 14.70565438827407
$ws.Cells.Item(22, 3).Value = 5.96278213881091
$ws.Cells.Item(22, 4).Value = 9.160802003894569
$ws.Cells.Item(22, 5).Value = 10.38059099726841
$ws.Cells.Item(22, 6).Value = 47.33775756086188
$ws.Cells.Item(22, 8).Value = 7.344005520526261
$ws.Cells.Item(22, 9).Value = 38.05864871053932
$ws.Cells.Item(22, 11).Value = 13.02003815903161
$ws.Cells.Item(22, 12).Value = 10.53922819678812
$ws.Cells.Item(22, 13).Value = 15.47828771314752
$ws.Cells.Item(23, 2).Value = 14.65431356334498
$ws.Cells.Item(23, 3).Value = 5.872242737331114
$ws.Cells.Item(23, 4).Value = 9.15656712112248
$ws.Cells.Item(23, 5).Value = 10.3836463740069
$ws.Cells.Item(23, 6).Value = 47.16279314071489
$ws.Cells.Item(23, 8).Value = 7.344005520526261
$ws.Cells.Item(23, 9).Value = 37.96912815309024
$ws.Cells.Item(23, 11).Value = 12.97998372391842
$ws.Cells.Item(23, 12).Value = 10.53820240897063
$ws.Cells.Item(23, 13).Value = 15.47365080236726
$ws.Cells.Item(24, 2).Value = 14.46641553655842
$ws.Cells.Item(24, 3).Value = 5.51455999228208
$ws.Cells.Item(24, 4).Value = 9.140082750533036
$ws.Cells.Item(24, 5).Value = 10.39658131870129
$ws.Cells.Item(24, 6).Value = 46.50530610530308
$ws.Cells.Item(24, 8).Value = 7.344005520526261
$ws.Cells.Item(24, 9).Value = 37.63497054660531
$ws.Cells.Item(24, 11).Value = 12.83400861539092
$ws.Cells.Item(24, 12).Value = 10.53687959137258
$ws.Cells.Item(24, 13).Value = 15.46144256690058
$ws.Cells.Item(25, 2).Value = 14.27875025774938
$ws.Cells.Item(25, 3).Value = 5.09748333206925
$ws.Cells.Item(25, 4).Value = 9.121340682533857
$ws.Cells.Item(25, 5).Value = 10.41345687919252
$ws.Cells.Item(25, 6).Value = 45.80954633866832
$ws.Cells.Item(25, 8).Value = 7.344005520526261
$ws.Cells.Item(25, 9).Value = 37.28581523774007
$ws.Cells.Item(25, 11).Value = 12.68951132431414
$ws.Cells.Item(25, 12).Value = 10.54094953229528
$ws.Cells.Item(25, 13).Value = 15.45977802025102
